$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (row 1) - match the style of the existing headers (e.g. H1)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

# New data cells (row 2)
$ws.Range("I2").Value = 9
$ws.Range("J2").Value = 9
